# Adding a qst to the dataset
#
# This Q&A workbook (sheet "Feuil1") stores a simple chat-bot script in
# columns A (question) / B (answer). The edit:
#   - fixes the "girle" -> "girl" typo in a few Q&A pairs
#   - fixes the developers' names ("ECCHARAY Mohammed" / "lahcen abdessalam"
#     -> "ECH-CHARAY Mohamed" / "IDALI LAHCEN Abdessalam") and their
#     casing/spelling
#   - adds a brand-new question/answer pair ("are you a boy or a girl" /
#     "what do you think ?") right after the "ok ! ?" / "good boy" row,
#     which shifts the rest of that little conversation down by one row
#     through row 143.
#
# Rows 1-131 and the B column of rows 1-18/20-123 are untouched. Only the
# A/B cell text below is rewritten to match the corrected conversation
# order; the surrounding styles/row heights already fit (no rows need to
# be inserted or deleted - row 143 stays the last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A132").Value = "are you a boy or a girl"
$ws.Range("B132").Value = "what do you think ?"

$ws.Range("A133").Value = "are you girl ?"
$ws.Range("B133").Value = "yes i m a girl !"

$ws.Range("A134").Value = "are you a boy ?"
$ws.Range("B134").Value = "no i m not like you i am girl ,seems you did not try using microphone ?"

$ws.Range("A135").Value = "how to use microphone ?"
$ws.Range("B135").Value = "juste click on the microphone icon to start recording when you end press it again Ok?"

$ws.Range("A136").Value = "are you a robot"
$ws.Range("B136").Value = "No, im not, i am only a bunch of codes deployed on a server "

$ws.Range("A137").Value = "whoes your developpers"
$ws.Range("B137").Value = "i'm developped by a team of two enginners students IDALI LAHCEN Abdessalam And ECH-CHARAY Mohamed from Ensah School"

$ws.Range("A138").Value = "tell me about yourself ?"
$ws.Range("B138").Value = "you know python ,flutter,deep learning i'm a mixt of those things what about you, tell me about your self what do you love ?"

$ws.Range("A139").Value = "guess what is my name ?"
$ws.Range("B139").Value = "really i have a very low memory ? Maybe i can keep your name in the next version ,Sorry"

$ws.Range("A140").Value = "nice to meet you"
$ws.Range("B140").Value = "nice to meet you too "

$ws.Range("A141").Value = "my name is"
$ws.Range("B141").Value = "nice to meet you"

$ws.Range("A142").Value = "tell me a joke"
$ws.Range("B142").Value = "i m not supposed to "

$ws.Range("A143").Value = "are IDALI LAHCEN Abdessalam And ECH-CHARAY Mohamed your developers?"
$ws.Range("B143").Value = "Yes of course, They did a great job. I am not a perfect bot but a great one. `nIf you want to test me, try asking me what ever you want about covid-19. Or simplly start a regular discussion with me."

# Update the saved view/selection to match (scrolled so row 129 is at the
# top, active cell on the new last question in column A).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 129
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A143").Select()
